# Revamped, offer section added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old data rows (rows 2-4) and rebuild with the single offer row.
# Row 1 headers stay the same (roll_no, company_name, package, role).
$ws.Range("A2:D4").ClearContents()

$ws.Range("A2").Value = "19IT1058"
$ws.Range("B2").Value = "TruckWawale"
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = "Data Analyst"

# Remove now-empty rows 3 and 4 so the used range shrinks back to A1:D2
$ws.Rows("3:4").Delete()

# Update the selected cell to match the post-edit workbook state
$ws.Range("B4").Select()
